# Sprint_1.xlsx edit — "Dopisałem się do sprinta. K Ziel"
#
# Adds Jakub Zielinski's (GUI) row-4 sprint entries, widens the
# "Planowane zadania" / "Przewidywany czas" columns so the longer text
# fits, and leaves the selection on the newly-filled cell D15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Jakub Zielinski / GUI) sprint log -----------------------------
# Written in this exact order so the new shared-string table entries come
# out in the same sequence as the target workbook (96=plan, 97=done,
# 98=time spent, 99=time expected).
$ws.Range("D4").Value = "Diagramy klas, ogarnięcie pojawiajacych się problemów z postawieniem projektu, opieprzanie Pauliny za inicjatywę godną rzodkiewki."
$ws.Range("F4").Value = "postawiłem projekt!"
$ws.Range("G4").Value = "8h"
$ws.Range("E4").Value = "3/4dni"

# --- Widen columns D & E so the new text is readable ----------------------
$ws.Columns("D").ColumnWidth = 122.5
$ws.Columns("E").ColumnWidth = 24.5

# --- Leave the selection where the edit happened ---------------------------
$ws.Range("D15").Select() | Out-Null
